# Trading update: 2026-02-18 10:25:59
$wb = $excel.ActiveWorkbook

# Helper: write a plain-text date string (e.g. "2026-02-18") into a cell
# without Excel's autodetection turning it into a date serial number.
function Set-TextDate($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# ---------------------------------------------------------------------
# 1. Add two new worksheets at the end of the workbook: "momentum" and
#    "MarketMaking" (per-strategy filtered trade logs).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMomentum = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMomentum.Name = "momentum"
$wsMomentum.PageSetup.LeftMargin = 0.75 * 72
$wsMomentum.PageSetup.RightMargin = 0.75 * 72
$wsMomentum.PageSetup.TopMargin = 1 * 72
$wsMomentum.PageSetup.BottomMargin = 1 * 72
$wsMomentum.PageSetup.HeaderMargin = 0.5 * 72
$wsMomentum.PageSetup.FooterMargin = 0.5 * 72

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMarketMaking = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMarketMaking.Name = "MarketMaking"
$wsMarketMaking.PageSetup.LeftMargin = 0.75 * 72
$wsMarketMaking.PageSetup.RightMargin = 0.75 * 72
$wsMarketMaking.PageSetup.TopMargin = 1 * 72
$wsMarketMaking.PageSetup.BottomMargin = 1 * 72
$wsMarketMaking.PageSetup.HeaderMargin = 0.5 * 72
$wsMarketMaking.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------
# 2. Summary sheet updates
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 1500   # Initial Capital
$wsSummary.Range("B3").Value = 1500   # Current Capital
$wsSummary.Range("B11").Value = 15    # Active Strategies

# ---------------------------------------------------------------------
# 3. Strategy Status sheet: populate the full strategy roster
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")

$strategies = @(
    "EMAArbitrage",
    "HighProbConvergence",
    "HighProbabilityBond",
    "LongshotBias",
    "MarketMaking",
    "MicrostructureScalper",
    "arbitrage",
    "breakout_momentum",
    "leadlag",
    "momentum",
    "orderbook_imbalance",
    "sentiment",
    "sharp_money",
    "volatility_scorer",
    "vwap"
)

$row = 2
foreach ($s in $strategies) {
    $wsStatus.Cells.Item($row, 1).Value = $s        # Strategy
    $wsStatus.Cells.Item($row, 2).Value = "ACTIVE"  # Status
    $wsStatus.Cells.Item($row, 3).Value = 100        # Capital
    $wsStatus.Cells.Item($row, 4).Value = 0          # Trades
    $wsStatus.Cells.Item($row, 5).Value = 0          # P&L $
    $wsStatus.Cells.Item($row, 6).Value = 0          # P&L %
    $wsStatus.Cells.Item($row, 7).Value = 0          # Win Rate %
    $row++
}

# ---------------------------------------------------------------------
# 4. All Trades sheet: add new columns K:Q and append new open trades
# ---------------------------------------------------------------------
$wsTrades = $wb.Worksheets.Item("All Trades")

# 4a. New header cells K1:Q1
$wsTrades.Range("K1").Value = "Capital After"
$wsTrades.Range("L1").Value = "Entry Slippage (bps)"
$wsTrades.Range("M1").Value = "Exit Slippage (bps)"
$wsTrades.Range("N1").Value = "Confidence"
$wsTrades.Range("O1").Value = "Entry Reason"
$wsTrades.Range("P1").Value = "Exit Reason"
$wsTrades.Range("Q1").Value = "Duration (min)"

# 4b. New trades appended as rows 19-22
$newTrades = @(
    @{ Row=19; Num=18; Date="2026-02-18"; Time="10:24:15"; Strategy="MarketMaking"; Side="UP";   Entry=0.67;     Status="OPEN"; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; EntryReason="Normal spread capture: 198 bps"; Duration=0 },
    @{ Row=20; Num=19; Date="2026-02-18"; Time="10:24:27"; Strategy="MarketMaking"; Side="DOWN"; Entry=0.428169; Status="OPEN"; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; EntryReason="Normal spread capture: 198 bps"; Duration=0 },
    @{ Row=21; Num=20; Date="2026-02-18"; Time="10:24:34"; Strategy="MarketMaking"; Side="DOWN"; Entry=0.48;     Status="OPEN"; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.6; EntryReason="Normal spread capture: 198 bps"; Duration=0 },
    @{ Row=22; Num=21; Date="2026-02-18"; Time="10:25:10"; Strategy="momentum";     Side="DOWN"; Entry=0.29;     Status="OPEN"; Capital=100; EntrySlip=0; ExitSlip=0; Confidence=0.9; EntryReason="Downward momentum: -1.980% over 10 samples"; Duration=0 }
)

foreach ($t in $newTrades) {
    $r = $t.Row
    $wsTrades.Cells.Item($r, 1).Value = $t.Num          # Trade #
    Set-TextDate $wsTrades.Cells.Item($r, 2) $t.Date    # Date
    $wsTrades.Cells.Item($r, 3).Value = $t.Time         # Time
    $wsTrades.Cells.Item($r, 4).Value = $t.Strategy     # Strategy
    $wsTrades.Cells.Item($r, 5).Value = $t.Side         # Side
    $wsTrades.Cells.Item($r, 6).Value = $t.Entry        # Entry Price
    # Exit Price (G) stays blank -- trade is still open
    $wsTrades.Cells.Item($r, 8).Value = $t.Status       # Status
    $wsTrades.Cells.Item($r, 9).Value = 0                # P&L %
    $wsTrades.Cells.Item($r, 10).Value = 0               # P&L $
    $wsTrades.Cells.Item($r, 11).Value = $t.Capital     # Capital After
    $wsTrades.Cells.Item($r, 12).Value = $t.EntrySlip   # Entry Slippage (bps)
    $wsTrades.Cells.Item($r, 13).Value = $t.ExitSlip    # Exit Slippage (bps)
    $wsTrades.Cells.Item($r, 14).Value = $t.Confidence  # Confidence
    $wsTrades.Cells.Item($r, 15).Value = $t.EntryReason # Entry Reason
    # Exit Reason (P) stays blank -- trade is still open
    $wsTrades.Cells.Item($r, 17).Value = $t.Duration    # Duration (min)
}

# ---------------------------------------------------------------------
# 5. "momentum" sheet: header + filtered momentum trade
# ---------------------------------------------------------------------
$headerCols = @("Trade #","Date","Time","Strategy","Side","Entry Price","Exit Price","Status","P&L %","P&L $","Capital After","Entry Slippage (bps)","Exit Slippage (bps)","Confidence","Entry Reason","Exit Reason","Duration (min)")
for ($c = 1; $c -le $headerCols.Length; $c++) {
    $wsMomentum.Cells.Item(1, $c).Value = $headerCols[$c - 1]
}

$wsMomentum.Cells.Item(2, 1).Value = 21
Set-TextDate $wsMomentum.Cells.Item(2, 2) "2026-02-18"
$wsMomentum.Cells.Item(2, 3).Value = "10:25:10"
$wsMomentum.Cells.Item(2, 4).Value = "momentum"
$wsMomentum.Cells.Item(2, 5).Value = "DOWN"
$wsMomentum.Cells.Item(2, 6).Value = 0.29
# Exit Price (G) stays blank
$wsMomentum.Cells.Item(2, 8).Value = "OPEN"
$wsMomentum.Cells.Item(2, 9).Value = 0
$wsMomentum.Cells.Item(2, 10).Value = 0
$wsMomentum.Cells.Item(2, 11).Value = 100
$wsMomentum.Cells.Item(2, 12).Value = 0
$wsMomentum.Cells.Item(2, 13).Value = 0
$wsMomentum.Cells.Item(2, 14).Value = 0.9
$wsMomentum.Cells.Item(2, 15).Value = "Downward momentum: -1.980% over 10 samples"
# Exit Reason (P) stays blank
$wsMomentum.Cells.Item(2, 17).Value = 0

# ---------------------------------------------------------------------
# 6. "MarketMaking" sheet: header + 3 filtered MarketMaking trades
# ---------------------------------------------------------------------
for ($c = 1; $c -le $headerCols.Length; $c++) {
    $wsMarketMaking.Cells.Item(1, $c).Value = $headerCols[$c - 1]
}

$mmTrades = @(
    @{ Row=2; Num=18; Date="2026-02-18"; Time="10:24:15"; Side="UP";   Entry=0.67;     Confidence=0.6 },
    @{ Row=3; Num=19; Date="2026-02-18"; Time="10:24:27"; Side="DOWN"; Entry=0.428169; Confidence=0.6 },
    @{ Row=4; Num=20; Date="2026-02-18"; Time="10:24:34"; Side="DOWN"; Entry=0.48;     Confidence=0.6 }
)

foreach ($t in $mmTrades) {
    $r = $t.Row
    $wsMarketMaking.Cells.Item($r, 1).Value = $t.Num
    Set-TextDate $wsMarketMaking.Cells.Item($r, 2) $t.Date
    $wsMarketMaking.Cells.Item($r, 3).Value = $t.Time
    $wsMarketMaking.Cells.Item($r, 4).Value = "MarketMaking"
    $wsMarketMaking.Cells.Item($r, 5).Value = $t.Side
    $wsMarketMaking.Cells.Item($r, 6).Value = $t.Entry
    # Exit Price (G) stays blank
    $wsMarketMaking.Cells.Item($r, 8).Value = "OPEN"
    $wsMarketMaking.Cells.Item($r, 9).Value = 0
    $wsMarketMaking.Cells.Item($r, 10).Value = 0
    $wsMarketMaking.Cells.Item($r, 11).Value = 100
    $wsMarketMaking.Cells.Item($r, 12).Value = 0
    $wsMarketMaking.Cells.Item($r, 13).Value = 0
    $wsMarketMaking.Cells.Item($r, 14).Value = $t.Confidence
    $wsMarketMaking.Cells.Item($r, 15).Value = "Normal spread capture: 198 bps"
    # Exit Reason (P) stays blank
    $wsMarketMaking.Cells.Item($r, 17).Value = 0
}

# ---------------------------------------------------------------------
# Keep the originally-active sheet ("Summary") selected, since the new
# sheets were appended after it and Excel would otherwise activate the
# most-recently-added sheet.
# ---------------------------------------------------------------------
$wsSummary.Activate()

Write-Output "Edit complete"
